# Auto-generated edit script applying cached market-price updates
# from the scheduled runner commit, per the supplied OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 992.2222
$ws.Range("I19").Value = 898.5
$ws.Range("J19").Value = 1067.2
$ws.Range("K19").Value = 898.5
$ws.Range("L19").Value = 1067.2
$ws.Range("M19").Value = -723.5
$ws.Range("N19").Value = -1417.2

$ws.Range("H33").Value = 270.24
$ws.Range("I33").Value = 120.44444
$ws.Range("J33").Value = 655.4286
$ws.Range("K33").Value = 120.44444
$ws.Range("L33").Value = 655.4286
$ws.Range("M33").Value = 108.55556
$ws.Range("N33").Value = -1113.4286

$ws.Range("H62").Value = 5055.148
$ws.Range("I62").Value = 4631.9
$ws.Range("K62").Value = 4631.9
$ws.Range("M62").Value = -4007.9

$ws.Range("H65").Value = 5055.148
$ws.Range("I65").Value = 4631.9
$ws.Range("K65").Value = 23159.5
$ws.Range("M65").Value = -20039.5

$ws.Range("H100").Value = 6448.5835
$ws.Range("I100").Value = 6987.6
$ws.Range("K100").Value = 6987.6
$ws.Range("M100").Value = -6446.6

$ws.Range("H116").Value = 20930.895
$ws.Range("I116").Value = 24421.643
$ws.Range("J116").Value = 11156.8
$ws.Range("K116").Value = 24421.643
$ws.Range("L116").Value = 11156.8
$ws.Range("M116").Value = -20979.643
$ws.Range("N116").Value = -18040.8

$ws.Range("H132").Value = 24576.58
$ws.Range("I132").Value = 27621.297
$ws.Range("K132").Value = 82863.891
$ws.Range("M132").Value = -80333.891

$ws.Range("H137").Value = 16587.3
$ws.Range("I137").Value = 21387.066
$ws.Range("J137").Value = 2188
$ws.Range("K137").Value = 64161.198
$ws.Range("L137").Value = 6564
$ws.Range("M137").Value = -61611.198
$ws.Range("N137").Value = -11664

$ws.Range("H138").Value = 32722.121
$ws.Range("J138").Value = 94370.82000000001
$ws.Range("L138").Value = 283112.46
$ws.Range("N138").Value = -293392.46

$ws.Range("H141").Value = 4404.4287
$ws.Range("I141").Value = 4605.3335
$ws.Range("K141").Value = 13816.0005
$ws.Range("M141").Value = -8636.000499999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20078.34
$ws.Range("I32").Value = 20406.77
$ws.Range("K32").Value = 20406.77
$ws.Range("M32").Value = -20119.77

$ws.Range("H61").Value = 5254.185
$ws.Range("I61").Value = 1018.4211
$ws.Range("K61").Value = 1018.4211
$ws.Range("M61").Value = -806.4211

$ws.Range("H74").Value = 408323.88
$ws.Range("I74").Value = 857663.9
$ws.Range("K74").Value = 857663.9
$ws.Range("M74").Value = -856789.9

$ws.Range("H77").Value = 408323.88
$ws.Range("I77").Value = 857663.9
$ws.Range("K77").Value = 4288319.5
$ws.Range("M77").Value = -4283951.5

$ws.Range("H122").Value = 2779.2222
$ws.Range("I122").Value = 2166.6667
$ws.Range("J122").Value = 4004.3333
$ws.Range("K122").Value = 6500.000100000001
$ws.Range("L122").Value = 12012.9999
$ws.Range("M122").Value = -4050.000100000001
$ws.Range("N122").Value = -16912.9999

$ws.Range("H136").Value = 5254.185
$ws.Range("I136").Value = 1018.4211
$ws.Range("K136").Value = 3055.2633
$ws.Range("M136").Value = -505.2633000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1977.1875
$ws.Range("I86").Value = 1619.4286
$ws.Range("J86").Value = 2255.4443
$ws.Range("K86").Value = 1619.4286
$ws.Range("L86").Value = 2255.4443
$ws.Range("M86").Value = -496.4286
$ws.Range("N86").Value = -4501.4443

$ws.Range("H89").Value = 1977.1875
$ws.Range("I89").Value = 1619.4286
$ws.Range("J89").Value = 2255.4443
$ws.Range("K89").Value = 8097.143
$ws.Range("L89").Value = 11277.2215
$ws.Range("M89").Value = -2481.143
$ws.Range("N89").Value = -22509.2215

$ws.Range("H99").Value = 1167.1333
$ws.Range("I99").Value = 1000.5833
$ws.Range("K99").Value = 1000.5833
$ws.Range("M99").Value = 497.4167

$ws.Range("H105").Value = 2485
$ws.Range("I105").Value = 1903.25
$ws.Range("K105").Value = 1903.25
$ws.Range("M105").Value = -156.25

$ws.Range("H134").Value = 4260.3
$ws.Range("I134").Value = 4115.1665
$ws.Range("K134").Value = 12345.4995
$ws.Range("M134").Value = -9810.499500000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4546711
$ws.Range("I31").Value = 5000882
$ws.Range("K31").Value = 5000882
$ws.Range("M31").Value = -5000587

$ws.Range("H34").Value = 4546711
$ws.Range("I34").Value = 5000882
$ws.Range("K34").Value = 5000882
$ws.Range("M34").Value = -5000680

$ws.Range("H58").Value = 1287.1875
$ws.Range("I58").Value = 1013
$ws.Range("K58").Value = 1013
$ws.Range("M58").Value = -810

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H99").Value = 3862.818
$ws.Range("I99").Value = 2099
$ws.Range("J99").Value = 6949.5
$ws.Range("K99").Value = 2099
$ws.Range("L99").Value = 6949.5
$ws.Range("M99").Value = -601
$ws.Range("N99").Value = -9945.5

$ws.Range("H122").Value = 2372.25
$ws.Range("I122").Value = 2746
$ws.Range("K122").Value = 8238
$ws.Range("M122").Value = -5788

$ws.Range("H126").Value = 3862.818
$ws.Range("I126").Value = 2099
$ws.Range("J126").Value = 6949.5
$ws.Range("K126").Value = 6297
$ws.Range("L126").Value = 20848.5
$ws.Range("M126").Value = -3827
$ws.Range("N126").Value = -25788.5

$ws.Range("H136").Value = 1287.1875
$ws.Range("I136").Value = 1013
$ws.Range("K136").Value = 3039
$ws.Range("M136").Value = -489

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 687.4286
$ws.Range("I5").Value = 703
$ws.Range("J5").Value = 648.5
$ws.Range("K5").Value = 2109
$ws.Range("L5").Value = 1945.5
$ws.Range("M5").Value = -1997
$ws.Range("N5").Value = -2169.5

$ws.Range("H135").Value = 687.4286
$ws.Range("I135").Value = 703
$ws.Range("J135").Value = 648.5
$ws.Range("K135").Value = 6327
$ws.Range("L135").Value = 5836.5
$ws.Range("M135").Value = -3792
$ws.Range("N135").Value = -10906.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 65279.25
$ws.Range("I102").Value = 73890.57000000001
$ws.Range("K102").Value = 73890.57000000001
$ws.Range("M102").Value = -72268.57000000001

$ws.Range("H132").Value = 2125.6365
$ws.Range("I132").Value = 1801.3334
$ws.Range("K132").Value = 5404.0002
$ws.Range("M132").Value = -2874.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3399.4285
$ws.Range("I40").Value = 2966
$ws.Range("K40").Value = 2966
$ws.Range("M40").Value = -2830

$ws.Range("H46").Value = 4751.364
$ws.Range("I46").Value = 1975.5
$ws.Range("J46").Value = 6337.5713
$ws.Range("K46").Value = 1975.5
$ws.Range("L46").Value = 6337.5713
$ws.Range("M46").Value = -1787.5
$ws.Range("N46").Value = -6713.5713

$ws.Range("H55").Value = 2825.6667
$ws.Range("I55").Value = 2714
$ws.Range("J55").Value = 2905.4285
$ws.Range("K55").Value = 2714
$ws.Range("L55").Value = 2905.4285
$ws.Range("M55").Value = -2541
$ws.Range("N55").Value = -3251.4285

$ws.Range("H82").Value = 1769.3572
$ws.Range("I82").Value = 1657.3
$ws.Range("J82").Value = 2049.5
$ws.Range("K82").Value = 1657.3
$ws.Range("L82").Value = 2049.5
$ws.Range("M82").Value = -1296.3
$ws.Range("N82").Value = -2771.5

$ws.Range("H85").Value = 1769.3572
$ws.Range("I85").Value = 1657.3
$ws.Range("J85").Value = 2049.5
$ws.Range("K85").Value = 1657.3
$ws.Range("L85").Value = 2049.5
$ws.Range("M85").Value = -409.3
$ws.Range("N85").Value = -4545.5

$ws.Range("H93").Value = 2667.125
$ws.Range("I93").Value = 1973
$ws.Range("J93").Value = 4749.5
$ws.Range("K93").Value = 1973
$ws.Range("L93").Value = 4749.5
$ws.Range("M93").Value = -725
$ws.Range("N93").Value = -7245.5

$ws.Range("H136").Value = 3755.0454
$ws.Range("I136").Value = 3246.8462
$ws.Range("K136").Value = 9740.5386
$ws.Range("M136").Value = -7190.5386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 33825.094
$ws.Range("I122").Value = 39011.305
$ws.Range("J122").Value = 2707.8333
$ws.Range("K122").Value = 117033.915
$ws.Range("L122").Value = 8123.499899999999
$ws.Range("M122").Value = -114583.915
$ws.Range("N122").Value = -13023.4999

$ws.Range("H132").Value = 2791.476
$ws.Range("I132").Value = 1889
$ws.Range("K132").Value = 5667
$ws.Range("M132").Value = -3137

$ws.Range("H136").Value = 16627.525
$ws.Range("I136").Value = 25034.479
$ws.Range("K136").Value = 75103.43700000001
$ws.Range("M136").Value = -72553.43700000001

